# Add a new worksheet "detail_lambda" after the last sheet (detail_V_dis)
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "detail_lambda"

# Populate header row + data rows
$ws.Cells.Item(1,1).Value = "exp"
$ws.Cells.Item(1,2).Value = "phi_0"
$ws.Cells.Item(1,3).Value = "dV_ges"
$ws.Cells.Item(1,4).Value = "eps_0"
$ws.Cells.Item(1,5).Value = "lambda"

$ws.Cells.Item(2,1).Value = "ye"
$ws.Cells.Item(2,2).Value = 0.00018
$ws.Cells.Item(2,3).Value = 160
$ws.Cells.Item(2,4).Value = 0.2
$ws.Cells.Item(2,5).Value = 0.99572192500000001

$ws.Cells.Item(3,1).Value = "ye"
$ws.Cells.Item(3,2).Value = 0.0002
$ws.Cells.Item(3,3).Value = 160
$ws.Cells.Item(3,4).Value = 0.2
$ws.Cells.Item(3,5).Value = 0.99572192500000001

$ws.Cells.Item(4,1).Value = "ye"
$ws.Cells.Item(4,2).Value = 0.00025
$ws.Cells.Item(4,3).Value = 160
$ws.Cells.Item(4,4).Value = 0.2
$ws.Cells.Item(4,5).Value = 0.99572192500000001

$ws.Cells.Item(5,1).Value = "ye"
$ws.Cells.Item(5,2).Value = 0.0003
$ws.Cells.Item(5,3).Value = 160
$ws.Cells.Item(5,4).Value = 0.2
$ws.Cells.Item(5,5).Value = 0.996256684

$ws.Cells.Item(6,1).Value = "ye"
$ws.Cells.Item(6,2).Value = 0.00035
$ws.Cells.Item(6,3).Value = 160
$ws.Cells.Item(6,4).Value = 0.2
$ws.Cells.Item(6,5).Value = 0.996256684

$ws.Cells.Item(7,1).Value = "ye"
$ws.Cells.Item(7,2).Value = 0.00018
$ws.Cells.Item(7,3).Value = 200
$ws.Cells.Item(7,4).Value = 0.2
$ws.Cells.Item(7,5).Value = 0.95080213899999999

$ws.Cells.Item(8,1).Value = "ye"
$ws.Cells.Item(8,2).Value = 0.0002
$ws.Cells.Item(8,3).Value = 200
$ws.Cells.Item(8,4).Value = 0.2
$ws.Cells.Item(8,5).Value = 0.98930481299999995

$ws.Cells.Item(9,1).Value = "ye"
$ws.Cells.Item(9,2).Value = 0.00025
$ws.Cells.Item(9,3).Value = 200
$ws.Cells.Item(9,4).Value = 0.2
$ws.Cells.Item(9,5).Value = 0.99358288800000005

$ws.Cells.Item(10,1).Value = "ye"
$ws.Cells.Item(10,2).Value = 0.0003
$ws.Cells.Item(10,3).Value = 200
$ws.Cells.Item(10,4).Value = 0.2
$ws.Cells.Item(10,5).Value = 0.99518716600000001

$ws.Cells.Item(11,1).Value = "ye"
$ws.Cells.Item(11,2).Value = 0.00035
$ws.Cells.Item(11,3).Value = 200
$ws.Cells.Item(11,4).Value = 0.2
$ws.Cells.Item(11,5).Value = 0.99732620299999997

$ws.Cells.Item(12,1).Value = "ye"
$ws.Cells.Item(12,2).Value = 0.00018
$ws.Cells.Item(12,3).Value = 200
$ws.Cells.Item(12,4).Value = 0.1
$ws.Cells.Item(12,5).Value = 0.96243386200000003

$ws.Cells.Item(13,1).Value = "ye"
$ws.Cells.Item(13,2).Value = 0.00018
$ws.Cells.Item(13,3).Value = 240
$ws.Cells.Item(13,4).Value = 0.2
$ws.Cells.Item(13,5).Value = 0.74973261999999996

$ws.Cells.Item(14,1).Value = "ye"
$ws.Cells.Item(14,2).Value = 0.0002
$ws.Cells.Item(14,3).Value = 240
$ws.Cells.Item(14,4).Value = 0.2
$ws.Cells.Item(14,5).Value = 0.97433155100000002

$ws.Cells.Item(15,1).Value = "ye"
$ws.Cells.Item(15,2).Value = 0.00025
$ws.Cells.Item(15,3).Value = 240
$ws.Cells.Item(15,4).Value = 0.2
$ws.Cells.Item(15,5).Value = 0.98609625700000003

$ws.Cells.Item(16,1).Value = "ye"
$ws.Cells.Item(16,2).Value = 0.0003
$ws.Cells.Item(16,3).Value = 240
$ws.Cells.Item(16,4).Value = 0.2
$ws.Cells.Item(16,5).Value = 0.99144385000000002

$ws.Cells.Item(17,1).Value = "ye"
$ws.Cells.Item(17,2).Value = 0.00035
$ws.Cells.Item(17,3).Value = 240
$ws.Cells.Item(17,4).Value = 0.2
$ws.Cells.Item(17,5).Value = 0.99518716600000001

$ws.Cells.Item(18,1).Value = "ye"
$ws.Cells.Item(18,2).Value = 0.00018
$ws.Cells.Item(18,3).Value = 240
$ws.Cells.Item(18,4).Value = 0.05
$ws.Cells.Item(18,5).Value = 0.93015873000000004

$ws.Cells.Item(19,1).Value = "ye"
$ws.Cells.Item(19,2).Value = 0.0002
$ws.Cells.Item(19,3).Value = 240
$ws.Cells.Item(19,4).Value = 0.05
$ws.Cells.Item(19,5).Value = 0.94973545000000004

$ws.Cells.Item(20,1).Value = "ye"
$ws.Cells.Item(20,2).Value = 0.00025
$ws.Cells.Item(20,3).Value = 240
$ws.Cells.Item(20,4).Value = 0.05
$ws.Cells.Item(20,5).Value = 0.97566137600000002

$ws.Cells.Item(21,1).Value = "ye"
$ws.Cells.Item(21,2).Value = 0.0003
$ws.Cells.Item(21,3).Value = 240
$ws.Cells.Item(21,4).Value = 0.05
$ws.Cells.Item(21,5).Value = 0.98783068799999996

$ws.Cells.Item(22,1).Value = "ye"
$ws.Cells.Item(22,2).Value = 0.00035
$ws.Cells.Item(22,3).Value = 240
$ws.Cells.Item(22,4).Value = 0.05
$ws.Cells.Item(22,5).Value = 0.99523809500000004

$ws.Cells.Item(23,1).Value = "ye"
$ws.Cells.Item(23,2).Value = 0.00018
$ws.Cells.Item(23,3).Value = 240
$ws.Cells.Item(23,4).Value = 0.1
$ws.Cells.Item(23,5).Value = 0.92751322800000002

$ws.Cells.Item(24,1).Value = "ye"
$ws.Cells.Item(24,2).Value = 0.00018
$ws.Cells.Item(24,3).Value = 280
$ws.Cells.Item(24,4).Value = 0.2
$ws.Cells.Item(24,5).Value = 0.64973261999999998

$ws.Cells.Item(25,1).Value = "ye"
$ws.Cells.Item(25,2).Value = 0.0002
$ws.Cells.Item(25,3).Value = 280
$ws.Cells.Item(25,4).Value = 0.2
$ws.Cells.Item(25,5).Value = 0.91657754000000002

$ws.Cells.Item(26,1).Value = "ye"
$ws.Cells.Item(26,2).Value = 0.00025
$ws.Cells.Item(26,3).Value = 280
$ws.Cells.Item(26,4).Value = 0.2
$ws.Cells.Item(26,5).Value = 0.97754010700000005

$ws.Cells.Item(27,1).Value = "ye"
$ws.Cells.Item(27,2).Value = 0.0003
$ws.Cells.Item(27,3).Value = 280
$ws.Cells.Item(27,4).Value = 0.2
$ws.Cells.Item(27,5).Value = 0.98716577500000002

$ws.Cells.Item(28,1).Value = "ye"
$ws.Cells.Item(28,2).Value = 0.00035
$ws.Cells.Item(28,3).Value = 280
$ws.Cells.Item(28,4).Value = 0.2
$ws.Cells.Item(28,5).Value = 0.99251336899999998

$ws.Cells.Item(29,1).Value = "ye"
$ws.Cells.Item(29,2).Value = 0.00018
$ws.Cells.Item(29,3).Value = 280
$ws.Cells.Item(29,4).Value = 0.05
$ws.Cells.Item(29,5).Value = 0.89470899500000001

$ws.Cells.Item(30,1).Value = "ye"
$ws.Cells.Item(30,2).Value = 0.0002
$ws.Cells.Item(30,3).Value = 280
$ws.Cells.Item(30,4).Value = 0.05
$ws.Cells.Item(30,5).Value = 0.91216931199999995

$ws.Cells.Item(31,1).Value = "ye"
$ws.Cells.Item(31,2).Value = 0.00025
$ws.Cells.Item(31,3).Value = 280
$ws.Cells.Item(31,4).Value = 0.05
$ws.Cells.Item(31,5).Value = 0.96137566100000005

$ws.Cells.Item(32,1).Value = "ye"
$ws.Cells.Item(32,2).Value = 0.0003
$ws.Cells.Item(32,3).Value = 280
$ws.Cells.Item(32,4).Value = 0.05
$ws.Cells.Item(32,5).Value = 0.98148148099999999

$ws.Cells.Item(33,1).Value = "ye"
$ws.Cells.Item(33,2).Value = 0.00035
$ws.Cells.Item(33,3).Value = 280
$ws.Cells.Item(33,4).Value = 0.05
$ws.Cells.Item(33,5).Value = 0.98994709000000003

$ws.Cells.Item(34,1).Value = "ye"
$ws.Cells.Item(34,2).Value = 0.00018
$ws.Cells.Item(34,3).Value = 280
$ws.Cells.Item(34,4).Value = 0.1
$ws.Cells.Item(34,5).Value = 0.89894179900000004

# Match the saved selection state from the authored workbook
$ws.Range("F27").Select() | Out-Null

Write-Output "detail_lambda sheet added"
